# Update the "RES installed" sheet's RES Installed [MW] values (column C, rows 2-6).
# This is the source data driving all the downstream RANDBETWEEN-based formulas
# (Pg, Winter/Summer S1/S2/S3 sheets) and the Main!B7 SUM formula.
$wb = $excel.ActiveWorkbook

$wsRes = $wb.Worksheets.Item("RES installed")
$wsRes.Range("C2").Value = 5
$wsRes.Range("C3").Value = 5
$wsRes.Range("C4").Value = 1.5
$wsRes.Range("C5").Value = 1.5
$wsRes.Range("C6").Value = 1.5

# Force a full recalculation so that all dependent formulas (including the
# volatile RANDBETWEEN-based ones) pick up the new RES installed capacities.
$excel.CalculateFullRebuild()
